# Automatische test-sync: 2025-06-22 21:39:50
# Adds a new "Verzoek om factuur" log entry (row 43) to the "Logs" sheet and
# updates the "Dashboard" category breakdown to reflect it.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs: append row 43 -----------------------------------------------
$logs.Range("A43").Value = "Verzoek om factuur"
$logs.Range("B43").Value = "mailmind.test@zohomail.eu"
$logs.Range("C43").Value = "Kunt u mij een factuur sturen voor mijn laatste bestelling?"
$logs.Range("D43").Value = "Factuur / Administratie"

$antwoord43 = @"
Beste klant,
Bedankt voor je e-mail. Om je te helpen met de factuur voor je laatste bestelling, heb ik wat meer informatie nodig. Zou je mij alsjeblieft de volgende gegevens kunnen doorgeven: ordernummer en/of datum van de bestelling? Zodra ik deze gegevens heb ontvangen, zal ik ervoor zorgen dat de factuur naar je wordt verstuurd.
Met vriendelijke groet,
[Naam] E-mailassistent
"@
$logs.Range("E43").Value = $antwoord43

$logs.Range("F43").Value = "2025-06-22 21:39:22"
$logs.Range("G43").Value = "Ja"

# Multi-line cell content makes the engine auto-pin an explicit row height;
# AutoFit clears that back to the sheet default (matches the source rows,
# e.g. row 40/41, which also hold multi-line text but carry no explicit
# row height).
$logs.Rows.Item(43).AutoFit()

# Extend the conditional formatting ranges (D2:D42 -> D2:D43, G2:G42 -> G2:G43)
# while keeping each rule's priority / dxfId intact.
$dFmt = $logs.Range("D2:D42").FormatConditions.Item(1)
$dFmt.ModifyAppliesToRange($logs.Range("D2:D43"))

$gFmt = $logs.Range("G2:G42").FormatConditions.Item(1)
$gFmt.ModifyAppliesToRange($logs.Range("G2:G43"))

# --- Dashboard: re-rank the category counts -----------------------------
$dash.Range("A11").Value = "Factuur / Administratie"
$dash.Range("B11").Value = 2

$dash.Range("A12").Value = "Juridisch / Contract"
$dash.Range("B12").Value = 2

$dash.Range("A14").Value = "Uitnodiging / Evenement"
$dash.Range("B14").Value = 1
